# Add leaf-bmc to leaf connections and fix cmm ports

$wb = $excel.ActiveWorkbook

$wsInterSwitch  = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$wsNonCompute   = $wb.Worksheets.Item("NON_COMPUTE_NODES")
$wsHwMgmt       = $wb.Worksheets.Item("HARDWARE_MANAGEMENT")
$wsCompute      = $wb.Worksheets.Item("COMPUTE_NODES")

# ---------------------------------------------------------------------------
# INTER_SWITCH_LINKS: add two new rows (43, 44) describing leaf-bmc <-> leaf
# connections. Duplicate row 42's formatting first (so the new rows inherit
# the same style as the rest of the table), then overwrite the values.
# ---------------------------------------------------------------------------
$wsInterSwitch.Rows.Item(42).Copy()
$wsInterSwitch.Rows.Item(43).Insert()
$wsInterSwitch.Rows.Item(42).Copy()
$wsInterSwitch.Rows.Item(44).Insert()

$wsInterSwitch.Range("J43").Value = "sw-leaf-bmc-001"
$wsInterSwitch.Range("K43").Value = "x3000"
$wsInterSwitch.Range("L43").Value = "u37"
$wsInterSwitch.Range("O43").Value = 48
$wsInterSwitch.Range("P43").Value = "sw-leaf-001"
$wsInterSwitch.Range("Q43").Value = "x3000"
$wsInterSwitch.Range("R43").Value = "u39"
$wsInterSwitch.Range("T43").Value = 51

$wsInterSwitch.Range("J44").Value = "sw-leaf-bmc-001"
$wsInterSwitch.Range("K44").Value = "x3000"
$wsInterSwitch.Range("L44").Value = "u37"
$wsInterSwitch.Range("O44").Value = 47
$wsInterSwitch.Range("P44").Value = "sw-leaf-002"
$wsInterSwitch.Range("Q44").Value = "x3000"
$wsInterSwitch.Range("R44").Value = "u38"
$wsInterSwitch.Range("T44").Value = 51

# Column J got a bit wider to fit the new "sw-leaf-bmc-001" values
# (closest achievable value to the target 15.1640625 given COM rounding).
$wsInterSwitch.Columns.Item(10).ColumnWidth = 14.25

# ---------------------------------------------------------------------------
# HARDWARE_MANAGEMENT: the cmm ports actually land on u37, not u39 - fix
# column R for rows 15 through 24.
# ---------------------------------------------------------------------------
$wsHwMgmt.Range("R15").Value = "u37"
$wsHwMgmt.Range("R16").Value = "u37"
$wsHwMgmt.Range("R17").Value = "u37"
$wsHwMgmt.Range("R18").Value = "u37"
$wsHwMgmt.Range("R19").Value = "u37"
$wsHwMgmt.Range("R20").Value = "u37"
$wsHwMgmt.Range("R21").Value = "u37"
$wsHwMgmt.Range("R22").Value = "u37"
$wsHwMgmt.Range("R23").Value = "u37"
$wsHwMgmt.Range("R24").Value = "u37"

# ---------------------------------------------------------------------------
# COMPUTE_NODES: ports were off by one - bump O20:O23 from 1 to 2.
# ---------------------------------------------------------------------------
$wsCompute.Range("O20").Value = 2
$wsCompute.Range("O21").Value = 2
$wsCompute.Range("O22").Value = 2
$wsCompute.Range("O23").Value = 2

# ---------------------------------------------------------------------------
# View / selection bookkeeping to match where the author ended up.
# ---------------------------------------------------------------------------

# NON_COMPUTE_NODES: just scrolled a bit further down (selection unchanged).
$wsNonCompute.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 2

# HARDWARE_MANAGEMENT: scrolled right, and selected the fixed R column range.
$wsHwMgmt.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$wsHwMgmt.Range("R16:R24").Select()

# INTER_SWITCH_LINKS: scrolled down a bit further, selection moved to the
# newly-added columns, and it is no longer the active tab.
$wsInterSwitch.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 3
$wsInterSwitch.Range("O45").Select()

# COMPUTE_NODES becomes the active tab, selection on the mirrored column.
$wsCompute.Activate()
$wsCompute.Range("O23").Select()

$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 880
